$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 18.93154223749944
$ws.Cells.Item(2, 4).Value = 7.486252983760925
$ws.Cells.Item(2, 5).Value = 26.34420579266186
$ws.Cells.Item(2, 6).Value = 43.90995181043784
$ws.Cells.Item(2, 7).Value = 3.639151515505089
$ws.Cells.Item(2, 9).Value = 24.22231225039193
$ws.Cells.Item(2, 12).Value = 13.80559457550026
$ws.Cells.Item(3, 2).Value = 18.40400383739879
$ws.Cells.Item(3, 4).Value = 7.427461062668137
$ws.Cells.Item(3, 5).Value = 24.69422820864891
$ws.Cells.Item(3, 6).Value = 42.78218030129256
$ws.Cells.Item(3, 7).Value = 3.648051832361162
$ws.Cells.Item(3, 9).Value = 24.55071170708684
$ws.Cells.Item(3, 12).Value = 13.46098306542121
$ws.Cells.Item(4, 2).Value = 18.07535418724256
$ws.Cells.Item(4, 4).Value = 7.395980570209996
$ws.Cells.Item(4, 5).Value = 23.62263855435555
$ws.Cells.Item(4, 6).Value = 42.10184247013062
$ws.Cells.Item(4, 7).Value = 3.653768711543347
$ws.Cells.Item(4, 9).Value = 24.76092740271038
$ws.Cells.Item(4, 12).Value = 13.24884841513078
$ws.Cells.Item(5, 2).Value = 17.94042784500537
$ws.Cells.Item(5, 4).Value = 7.384309345525342
$ws.Cells.Item(5, 5).Value = 23.1712759940711
$ws.Cells.Item(5, 6).Value = 41.82808832876966
$ws.Cells.Item(5, 7).Value = 3.656162232874639
$ws.Cells.Item(5, 9).Value = 24.84876178970738
$ws.Cells.Item(5, 12).Value = 13.16239092067583
$ws.Cells.Item(6, 2).Value = 17.91796895185421
$ws.Cells.Item(6, 4).Value = 7.382441069519056
$ws.Cells.Item(6, 5).Value = 23.0954421209373
$ws.Cells.Item(6, 6).Value = 41.78285528461852
$ws.Cells.Item(6, 7).Value = 3.656563545024724
$ws.Cells.Item(6, 9).Value = 24.8634780516073
$ws.Cells.Item(6, 12).Value = 13.14803780794216
$ws.Cells.Item(7, 2).Value = 18.07353830496393
$ws.Cells.Item(7, 4).Value = 7.395818489543315
$ws.Cells.Item(7, 5).Value = 23.61661071206052
$ws.Cells.Item(7, 6).Value = 42.09813581555773
$ws.Cells.Item(7, 7).Value = 3.653800732279901
$ws.Cells.Item(7, 9).Value = 24.76210316536039
$ws.Cells.Item(7, 12).Value = 13.2476822944495
$ws.Cells.Item(8, 2).Value = 18.7507425284142
$ws.Cells.Item(8, 4).Value = 7.465021720837463
$ws.Cells.Item(8, 5).Value = 25.78744540290887
$ws.Cells.Item(8, 6).Value = 43.51884563279096
$ws.Cells.Item(8, 7).Value = 3.642168315410911
$ws.Cells.Item(8, 9).Value = 24.33377208995873
$ws.Cells.Item(8, 12).Value = 13.68695873563933
$ws.Cells.Item(9, 2).Value = 20.03328718501071
$ws.Cells.Item(9, 4).Value = 7.637385567601935
$ws.Cells.Item(9, 5).Value = 29.58002544393358
$ws.Cells.Item(9, 6).Value = 46.38181724292822
$ws.Cells.Item(9, 7).Value = 3.621335563657547
$ws.Cells.Item(9, 9).Value = 23.56127047856659
$ws.Cells.Item(9, 12).Value = 14.53883118912527
$ws.Cells.Item(10, 2).Value = 20.93826574673055
$ws.Cells.Item(10, 4).Value = 7.786174136098849
$ws.Cells.Item(10, 5).Value = 32.0840434632935
$ws.Cells.Item(10, 6).Value = 48.5073024520302
$ws.Cells.Item(10, 7).Value = 3.607205354309119
$ws.Cells.Item(10, 9).Value = 23.03398224515062
$ws.Cells.Item(10, 12).Value = 15.15228869094822
$ws.Cells.Item(11, 2).Value = 21.34009882768034
$ws.Cells.Item(11, 4).Value = 7.858570037809044
$ws.Cells.Item(11, 5).Value = 33.16225185366297
$ws.Cells.Item(11, 6).Value = 49.47413545366794
$ws.Cells.Item(11, 7).Value = 3.601025720509365
$ws.Cells.Item(11, 9).Value = 22.80266307115258
$ws.Cells.Item(11, 12).Value = 15.42741590334697
$ws.Cells.Item(12, 2).Value = 21.49071737719738
$ws.Cells.Item(12, 4).Value = 7.886649463528663
$ws.Cells.Item(12, 5).Value = 33.56182586500049
$ws.Cells.Item(12, 6).Value = 49.83986590875629
$ws.Cells.Item(12, 7).Value = 3.598720810316119
$ws.Cells.Item(12, 9).Value = 22.71628287263836
$ws.Cells.Item(12, 12).Value = 15.5309387906421
$ws.Cells.Item(13, 2).Value = 21.45834972722984
$ws.Cells.Item(13, 4).Value = 7.8805727225901
$ws.Cells.Item(13, 5).Value = 33.47615766231155
$ws.Cells.Item(13, 6).Value = 49.76112232265491
$ws.Cells.Item(13, 7).Value = 3.599215657547322
$ws.Cells.Item(13, 9).Value = 22.73483256006374
$ws.Cells.Item(13, 12).Value = 15.50867415697847
$ws.Cells.Item(14, 2).Value = 21.35252201387238
$ws.Cells.Item(14, 4).Value = 7.860866909222697
$ws.Cells.Item(14, 5).Value = 33.19529959061747
$ws.Cells.Item(14, 6).Value = 49.50423417329714
$ws.Cells.Item(14, 7).Value = 3.600835392046529
$ws.Cells.Item(14, 9).Value = 22.79553224169196
$ws.Cells.Item(14, 12).Value = 15.43594657681641
$ws.Cells.Item(15, 2).Value = 21.28749435991653
$ws.Cells.Item(15, 4).Value = 7.848882671362359
$ws.Cells.Item(15, 5).Value = 33.02213143487459
$ws.Cells.Item(15, 6).Value = 49.34682141894118
$ws.Cells.Item(15, 7).Value = 3.601832092516926
$ws.Cells.Item(15, 9).Value = 22.83287041572821
$ws.Cells.Item(15, 12).Value = 15.3913099438983
$ws.Cells.Item(16, 2).Value = 20.91179538171079
$ws.Cells.Item(16, 4).Value = 7.781536653118483
$ws.Cells.Item(16, 5).Value = 32.01235772059846
$ws.Cells.Item(16, 6).Value = 48.44408335064775
$ws.Cells.Item(16, 7).Value = 3.607614159688183
$ws.Cells.Item(16, 9).Value = 23.04927031601998
$ws.Cells.Item(16, 12).Value = 15.13422089356277
$ws.Cells.Item(17, 2).Value = 20.67869702014961
$ws.Cells.Item(17, 4).Value = 7.741420333980248
$ws.Cells.Item(17, 5).Value = 31.3773234542855
$ws.Cells.Item(17, 6).Value = 47.88998854715715
$ws.Cells.Item(17, 7).Value = 3.611224493041729
$ws.Cells.Item(17, 9).Value = 23.18420401634861
$ws.Cells.Item(17, 12).Value = 14.97542604592192
$ws.Cells.Item(18, 2).Value = 20.54370523478015
$ws.Cells.Item(18, 4).Value = 7.718790853988901
$ws.Cells.Item(18, 5).Value = 31.00633550637787
$ws.Cells.Item(18, 6).Value = 47.57130317205676
$ws.Cells.Item(18, 7).Value = 3.613324461141837
$ws.Cells.Item(18, 9).Value = 23.26261961155211
$ws.Cells.Item(18, 12).Value = 14.88372603486408
$ws.Cells.Item(19, 2).Value = 20.49784556902874
$ws.Cells.Item(19, 4).Value = 7.711205573782077
$ws.Cells.Item(19, 5).Value = 30.87973976038741
$ws.Cells.Item(19, 6).Value = 47.46341641886482
$ws.Cells.Item(19, 7).Value = 3.614039508817736
$ws.Cells.Item(19, 9).Value = 23.28930851035359
$ws.Cells.Item(19, 12).Value = 14.8526182679841
$ws.Cells.Item(20, 2).Value = 20.7036069219452
$ws.Cells.Item(20, 4).Value = 7.745644880065508
$ws.Cells.Item(20, 5).Value = 31.44551705271595
$ws.Cells.Item(20, 6).Value = 47.94897413794731
$ws.Cells.Item(20, 7).Value = 3.610837748817856
$ws.Cells.Item(20, 9).Value = 23.16975684684436
$ws.Cells.Item(20, 12).Value = 14.99236856502408
$ws.Cells.Item(21, 2).Value = 21.38364914206413
$ws.Cells.Item(21, 4).Value = 7.866637058322666
$ws.Cells.Item(21, 5).Value = 33.27803076807638
$ws.Cells.Item(21, 6).Value = 49.57970187198939
$ws.Cells.Item(21, 7).Value = 3.60035868607918
$ws.Cells.Item(21, 9).Value = 22.7776703991005
$ws.Cells.Item(21, 12).Value = 15.45732711413562
$ws.Cells.Item(22, 2).Value = 21.81902760974936
$ws.Cells.Item(22, 4).Value = 7.949579218606634
$ws.Cells.Item(22, 5).Value = 34.4248773116753
$ws.Cells.Item(22, 6).Value = 50.64305789105741
$ws.Cells.Item(22, 7).Value = 3.593714838662582
$ws.Cells.Item(22, 9).Value = 22.52849680158296
$ws.Cells.Item(22, 12).Value = 15.75730807138112
$ws.Cells.Item(23, 2).Value = 21.58752768374333
$ws.Cells.Item(23, 4).Value = 7.904962398654054
$ws.Cells.Item(23, 5).Value = 33.81741953736175
$ws.Cells.Item(23, 6).Value = 50.07586417825629
$ws.Cells.Item(23, 7).Value = 3.597242217238649
$ws.Cells.Item(23, 9).Value = 22.66084239658213
$ws.Cells.Item(23, 12).Value = 15.59758863889483
$ws.Cells.Item(24, 2).Value = 20.69234819908894
$ws.Cells.Item(24, 4).Value = 7.743733610250738
$ws.Cells.Item(24, 5).Value = 31.41470510127477
$ws.Cells.Item(24, 6).Value = 47.9223071247079
$ws.Cells.Item(24, 7).Value = 3.611012520030562
$ws.Cells.Item(24, 9).Value = 23.17628579344312
$ws.Cells.Item(24, 12).Value = 14.9847101116964
$ws.Cells.Item(25, 2).Value = 19.69223309491429
$ws.Cells.Item(25, 4).Value = 7.586840400505531
$ws.Cells.Item(25, 5).Value = 28.60372160413435
$ws.Cells.Item(25, 6).Value = 45.60169685859243
$ws.Cells.Item(25, 7).Value = 3.62676269494144
$ws.Cells.Item(25, 9).Value = 23.76311894796806
$ws.Cells.Item(25, 12).Value = 14.31010605028944
